$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5453776865001148
$ws.Range("E2").Value = -0.7976031984000098

$ws.Range("C3").Value = -6.170514117037273
$ws.Range("E3").Value = -8.396348489509153

$ws.Range("C4").Value = -0.2290082001396909
$ws.Range("E4").Value = -4.327930935900004

$ws.Range("C5").Value = 1.406827509327035
$ws.Range("E5").Value = 2.015050062499957

$ws.Range("C6").Value = 1.153683074671208
$ws.Range("E6").Value = 3.648892256099945

$ws.Range("C7").Value = 0.2186142574756467
$ws.Range("E7").Value = 0.4006004000999708

$ws.Range("C8").Value = -0.8522658067264599
$ws.Range("E8").Value = -3.551690943899999

$ws.Range("C9").Value = -0.2262139320475365
$ws.Range("E9").Value = -0.7976031983999876

$ws.Range("C10").Value = -0.6258176826215101
$ws.Range("E10").Value = -0.3994003999000073

$ws.Range("C11").Value = 0.3239252862367037
$ws.Range("E11").Value = 1.609625625600009

$ws.Range("C12").Value = 0.5738128002843901
$ws.Range("E12").Value = -0.3994003999000184

$ws.Range("C13").Value = -0.4781004700720293
$ws.Range("E13").Value = 0.8024032015999882

$ws.Range("C14").Value = -1.197849743493773
$ws.Range("E14").Value = -3.161804390399992

$ws.Range("C15").Value = 1.064698711638945
$ws.Range("E15").Value = -2.540956581357878

$ws.Range("C16").Value = -1.524103236349472
$ws.Range("E16").Value = -1.240907591477092

$ws.Range("C17").Value = 0.2854413827033664
$ws.Range("E17").Value = -0.2470349027347551

$ws.Range("C18").Value = 1.118108578853261
$ws.Range("E18").Value = 1.532721825047534

$ws.Range("C19").Value = -1.490505436658163
$ws.Range("E19").Value = -0.3349088112516219
